$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.384.24'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '1.868.21'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7035'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.40%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07921'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3133'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.48'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07814'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.96%  '
$ws.Range('D12').Value = '1.861.83'
$ws.Range('E12').Value = '  -0.69%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '93.93'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.175'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7003'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.535'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.94%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008391'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.43%  '
$ws.Range('D18').Value = '29.353.41'
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '252.60'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.57%  '
$ws.Range('D20').Value = '2.111.82'
$ws.Range('E20').Value = '  -0.83%  '
$ws.Range('E21').Value = '  -1.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9996'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.642'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.55%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.000'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1554'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.003'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '161.62'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.85'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.503'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.315'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.253'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.23%  '
$ws.Range('E32').Value = '  +3.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05274'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.895'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.175'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7486'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.708'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.64%  '
$ws.Range('E38').Value = '  +0.18%  '
$ws.Range('D39').Value = '1.272.58'
$ws.Range('E39').Value = '  +0.86%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.771'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8910'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.20%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.001'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.61%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '108.57'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '71.01'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.97%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000129'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.64%  '
$ws.Range('E47').Value = '  -0.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.614'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.800'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.5180'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.84%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4302'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.11%  '
